# Insert a new data row at row 13 (pushing the existing rows 13-23 down to
# 14-24), then populate the new row 13 with the new Pomelo price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 13:23 down by one row to make room for the new record, then
# fill the freed-up row 13 with the new data.
$ws.Rows.Item(13).Insert(-4121)   # xlShiftDown = -4121

$ws.Range("A13").Value = 6
$ws.Range("B13").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 44574
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100102
$ws.Range("H13").Value = "Cítricos"
$ws.Range("I13").Value = 100102006
$ws.Range("J13").Value = "Pomelo"
$ws.Range("K13").Value = "Start Ruby"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 16
$ws.Range("N13").Value = 200000
$ws.Range("O13").Value = 200000
$ws.Range("P13").Value = 200000
$ws.Range("Q13").Value = "`$/bins (350 kilos)"
$ws.Range("R13").Value = "Región Metropolitana"
$ws.Range("S13").Value = 571
$ws.Range("T13").Value = 350
